$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.ClearContents()

$ws.Range("A2").Value = "Contact"
$ws.Range("A2").Font.FontStyle = "Bold"
$ws.Range("B2").Value = "Test User User"
$ws.Range("G2").Value = "Data of Estimate:"
$ws.Range("G2").Font.FontStyle = "Bold"
$ws.Range("L2").Value = "Sales Rep:"
$ws.Range("L2").Font.FontStyle = "Bold"
$ws.Range("A3").Value = "Firm:"
$ws.Range("A3").Font.FontStyle = "Bold"
$ws.Range("G3").Value = "Project Number:"
$ws.Range("G3").Font.FontStyle = "Bold"
$ws.Range("L3").Value = "Region:"
$ws.Range("L3").Font.FontStyle = "Bold"
$ws.Range("A4").Value = "Address:"
$ws.Range("A4").Font.FontStyle = "Bold"
$ws.Range("G4").Value = "Project Name:"
$ws.Range("G4").Font.FontStyle = "Bold"
$ws.Range("A5").Value = "Email"
$ws.Range("A5").Font.FontStyle = "Bold"
$ws.Range("B5").Value = "user@user.com"
$ws.Range("G5").Value = "Project Address:"
$ws.Range("G5").Font.FontStyle = "Bold"
$ws.Range("L5").Value = "Lead Time:"
$ws.Range("L5").Font.FontStyle = "Bold"
$ws.Range("M5").Value = "8-10 WEEKS UPON ORDER APPROVAL AND RECEIPT OF DEPOSIT"
$ws.Range("A6").Value = "Phone Number"
$ws.Range("A6").Font.FontStyle = "Bold"
$ws.Range("A8").Value = "ROOM"
$ws.Range("A8").Font.FontStyle = "Bold"
$ws.Range("B8").Value = "REF #"
$ws.Range("B8").Font.FontStyle = "Bold"
$ws.Range("C8").Value = "QTY"
$ws.Range("C8").Font.FontStyle = "Bold"
$ws.Range("D8").Value = "FINISH"
$ws.Range("D8").Font.FontStyle = "Bold"
$ws.Range("E8").Value = "SIZE"
$ws.Range("E8").Font.FontStyle = "Bold"
$ws.Range("F8").Value = "PLATE COST"
$ws.Range("F8").Font.FontStyle = "Bold"
$ws.Range("G8").Value = "MECHANISM TYPE/QTY"
$ws.Range("G8").Font.FontStyle = "Bold"
$ws.Range("H8").Value = "MECHANISM PROVIDED BY"
$ws.Range("H8").Font.FontStyle = "Bold"
$ws.Range("I8").Value = "V&VERSER COST"
$ws.Range("I8").Font.FontStyle = "Bold"
$ws.Range("J8").Value = "TOTAL # ENGRAVING"
$ws.Range("J8").Font.FontStyle = "Bold"
$ws.Range("K8").Value = "ENGRAVING COST"
$ws.Range("K8").Font.FontStyle = "Bold"
$ws.Range("L8").Value = "BACK BOX"
$ws.Range("L8").Font.FontStyle = "Bold"
$ws.Range("M8").Value = "MELKIT"
$ws.Range("M8").Font.FontStyle = "Bold"
$ws.Range("N8").Value = "TOTAL"
$ws.Range("N8").Font.FontStyle = "Bold"
$ws.Range("O8").Value = "EDGES: STRAIGHT / BEVELED"
$ws.Range("O8").Font.FontStyle = "Bold"
$ws.Range("B9").Value = "FC3008A-A1100141FA"
$ws.Range("F9").Value = "$"
$ws.Range("H9").Value = "MELJAC"
$ws.Range("I9").Value = "$"
$ws.Range("K9").Value = "$"
$ws.Range("L9").Value = "$"
$ws.Range("M9").Value = "$"
$ws.Range("N9").Value = "$"
$ws.Range("O9").Value = "TBD"
$ws.Range("C9").Value = 1
$ws.Range("J9").Value = 0

$ws.Range("I2").WrapText = $false
$ws.Range("M2").WrapText = $false
$ws.Range("B3").WrapText = $false
$ws.Range("I3").WrapText = $false
$ws.Range("M3").WrapText = $false
$ws.Range("B4").WrapText = $false
$ws.Range("I4").WrapText = $false
$ws.Range("I5").WrapText = $false
$ws.Range("B6").WrapText = $false
$ws.Range("A9").WrapText = $false
$ws.Range("D9").WrapText = $false
$ws.Range("E9").WrapText = $false
$ws.Range("G9").WrapText = $false

$null = $ws.Range("O8").Select()
